$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44271
$ws.Range("J3").Value = 55
$ws.Range("M3").Value = 9227
$ws.Range("P3").Value = 154

# Row 4
$ws.Range("D4").Value = 44617
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 11000
$ws.Range("M4").Value = 10500
$ws.Range("O4").Value = 'Región Metropolitana'
$ws.Range("P4").Value = 175

# Row 5
$ws.Range("D5").Value = 44600
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12500
$ws.Range("O5").Value = 'Región de Arica y Parinacota'
$ws.Range("P5").Value = 208

# Row 6
$ws.Range("D6").Value = 44594
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("O6").Value = 'Región de Arica y Parinacota'
$ws.Range("P6").Value = 208

# Row 7
$ws.Range("D7").Value = 44160
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 7500
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = 7667
$ws.Range("P7").Value = 128

# Row 8
$ws.Range("D8").Value = 44224
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 8500
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 8719
$ws.Range("O8").Value = 'Región del Maule'
$ws.Range("P8").Value = 145

# Row 9
$ws.Range("D9").Value = 44259
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 9500
$ws.Range("M9").Value = 9214
$ws.Range("P9").Value = 154

# Row 10
$ws.Range("D10").Value = 44204
$ws.Range("J10").Value = 45
$ws.Range("K10").Value = 9500
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 9722
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 162

# Row 11
$ws.Range("D11").Value = 44159
$ws.Range("J11").Value = 35
$ws.Range("K11").Value = 7500
$ws.Range("L11").Value = 8000
$ws.Range("M11").Value = 7714
$ws.Range("O11").Value = 'Región de Arica y Parinacota'
$ws.Range("P11").Value = 129

# Row 12
$ws.Range("D12").Value = 44266
$ws.Range("J12").Value = 60
$ws.Range("K12").Value = 9000
$ws.Range("L12").Value = 9500
$ws.Range("M12").Value = 9208
$ws.Range("O12").Value = 'Región del Maule'
$ws.Range("P12").Value = 153

# Row 13
$ws.Range("D13").Value = 44253
$ws.Range("J13").Value = 95
$ws.Range("K13").Value = 9500
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 9658
$ws.Range("O13").Value = 'Región del Maule'
$ws.Range("P13").Value = 161

# Row 14
$ws.Range("D14").Value = 44216
$ws.Range("J14").Value = 55
$ws.Range("K14").Value = 9500
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = 9773
$ws.Range("P14").Value = 163

# Row 15
$ws.Range("D15").Value = 44210
$ws.Range("J15").Value = 60
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 8417
$ws.Range("O15").Value = 'Región de Arica y Parinacota'
$ws.Range("P15").Value = 140

# Row 16
$ws.Range("D16").Value = 44208
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = 7350
$ws.Range("O16").Value = 'Región del Maule'
$ws.Range("P16").Value = 122

# Row 17
$ws.Range("D17").Value = 44610
$ws.Range("J17").Value = 60
$ws.Range("K17").Value = 11000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 11500
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 192

# Row 18
$ws.Range("D18").Value = 44218
$ws.Range("J18").Value = 65
$ws.Range("L18").Value = 10000
$ws.Range("M18").Value = 9615
$ws.Range("P18").Value = 160

# Row 19
$ws.Range("D19").Value = 44202
$ws.Range("J19").Value = 50
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = 8400
$ws.Range("O19").Value = 'Región del Maule'
$ws.Range("P19").Value = 140

# Row 21
$ws.Range("D21").Value = 44596
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 13000
$ws.Range("M21").Value = 12500
$ws.Range("O21").Value = 'Región de Arica y Parinacota'
$ws.Range("P21").Value = 208

# Row 22
$ws.Range("D22").Value = 44608
$ws.Range("K22").Value = 12000
$ws.Range("L22").Value = 13000
$ws.Range("M22").Value = 12500
$ws.Range("P22").Value = 208

# Row 23
$ws.Range("D23").Value = 44624
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = 10000
$ws.Range("L23").Value = 11000
$ws.Range("M23").Value = 10500
$ws.Range("O23").Value = 'Región Metropolitana'
$ws.Range("P23").Value = 175

# Row 24
$ws.Range("D24").Value = 44162
$ws.Range("J24").Value = 43
$ws.Range("K24").Value = 8000
$ws.Range("L24").Value = 8500
$ws.Range("M24").Value = 8209
$ws.Range("O24").Value = 'Región de Arica y Parinacota'
$ws.Range("P24").Value = 137

# Row 26
$ws.Range("D26").Value = 44264
$ws.Range("J26").Value = 43
$ws.Range("K26").Value = 8500
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = 8709
$ws.Range("O26").Value = 'Región del Maule'
$ws.Range("P26").Value = 145
